$wb = $excel.ActiveWorkbook

# --- Overview sheet: status columns (zh-cn / de-de) reflect the new "handed back" status ---
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E2").Value = "Handed back: in sync with en-US"
$wsOverview.Range("F2").Value = "Handed back: in sync with en-US"

# --- zh-cn sheet: status, handback datetime and cleared error detail ---
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C2").Value = "Handed back: in sync with en-US"
$wsZhCn.Range("K2").Value = "2016-08-23 06:47:13"
$wsZhCn.Range("P2").Value = ""

# --- de-de sheet: status, handback datetime and cleared error detail ---
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C2").Value = "Handed back: in sync with en-US"
$wsDeDe.Range("K2").Value = "2016-08-23 06:47:20"
$wsDeDe.Range("P2").Value = ""

# --- Column width adjustments (report generation widens the Status columns and
#     shrinks the now largely-empty Error Detail column). Excel snaps
#     ColumnWidth to whole-pixel boundaries, so these inputs are chosen to
#     land on the closest achievable stored width. ---
$wsOverview.Columns.Item(5).ColumnWidth = 29.166666666666668
$wsOverview.Columns.Item(6).ColumnWidth = 29.166666666666668

$wsZhCn.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsZhCn.Columns.Item(16).ColumnWidth = 12.833333333333332

$wsDeDe.Columns.Item(3).ColumnWidth = 29.166666666666668
$wsDeDe.Columns.Item(16).ColumnWidth = 12.833333333333332
